$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ASV_rank (column G) for rows 9-12 from 41 to 42
$ws.Range("G9:G12").Value = 42
